$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# Replace the "==" separator runs with "*"
$find.Execute("==", $false, $false, $false, $false, $false, $true, 1, $false, "*", 2)
